# Auto-generated Excel COM-interop script
# Applies scheduled market-price / profit-data updates to the Leve profit tables
# across all crafting class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# For each affected Leve row, columns H-N are refreshed:
#   H = currentAveragePrice        K = LevePriceNQ
#   I = currentAveragePriceNQ      L = LevePriceHQ
#   J = currentAveragePriceHQ      M = LeveProfitNQ (present only when K <> 0)
#                                  N = LeveProfitHQ (present only when L <> 0)

$wb = $excel.ActiveWorkbook

# ---------------- ALC ----------------
$ws = $wb.Worksheets.Item("ALC")

# Row 57: Quit Your Jib-jab | Gold Needle
$ws.Range("H57").Value = 23000
$ws.Range("J57").Value = 23000
$ws.Range("L57").Value = 69000
$ws.Range("N57").Value = -69998

# Row 88: The Grave of Hemlock Groves | Growth Formula Zeta
$ws.Range("H88").Value = 24816.727
$ws.Range("I88").Value = 9333.333000000001
$ws.Range("J88").Value = 43396.8
$ws.Range("K88").Value = 9333.333000000001
$ws.Range("L88").Value = 43396.8
$ws.Range("M88").Value = -8927.333000000001
$ws.Range("N88").Value = -44208.8

# Row 91: Dappling the Highlands (L) | Growth Formula Zeta
$ws.Range("H91").Value = 24816.727
$ws.Range("I91").Value = 9333.333000000001
$ws.Range("J91").Value = 43396.8
$ws.Range("K91").Value = 9333.333000000001
$ws.Range("L91").Value = 43396.8
$ws.Range("M91").Value = -7929.333000000001
$ws.Range("N91").Value = -46204.8

# Row 108: Keeping Magic Alive | Smilodonskin Grimoire
$ws.Range("H108").Value = 33247
$ws.Range("J108").Value = 33247
$ws.Range("L108").Value = 33247
$ws.Range("N108").Value = -40927

# Row 117: A Greater Grimoire | Zonureskin Grimoire
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# Row 120: Supreme Official Strategy Guide | Dwarven Mythril Codex
$ws.Range("H120").Value = 48992
$ws.Range("J120").Value = 48992
$ws.Range("L120").Value = 48992
$ws.Range("N120").Value = -58668

# Row 129: Practical Command | Commanding Craftsman's Draught
$ws.Range("H129").Value = 232151.58
$ws.Range("I129").Value = 6005098.5
$ws.Range("J129").Value = 1233.7
$ws.Range("K129").Value = 18015295.5
$ws.Range("L129").Value = 3701.1
$ws.Range("M129").Value = -18010295.5
$ws.Range("N129").Value = -13701.1

# ---------------- ARM ----------------
$ws = $wb.Worksheets.Item("ARM")

# Row 23: A Well-rounded Crew | Iron Hoplon
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

# Row 37: Get Shirty | Steel Chainmail
$ws.Range("H37").Value = 32507.143
$ws.Range("I37").Value = 2400
$ws.Range("J37").Value = 37525
$ws.Range("K37").Value = 2400
$ws.Range("L37").Value = 37525
$ws.Range("M37").Value = -2127
$ws.Range("N37").Value = -38071

# Row 44: Very Slow Array | Mythril Plate
$ws.Range("H44").Value = 36375.23
$ws.Range("J44").Value = 36375.23
$ws.Range("L44").Value = 36375.23
$ws.Range("N44").Value = -37351.23

# Row 55: Employee Retention | Mythril Elmo
$ws.Range("H55").Value = 51000
$ws.Range("J55").Value = 51000
$ws.Range("L55").Value = 51000
$ws.Range("N55").Value = -51630

# Row 80: A Squire to Inspire | Titanium Hoplon
$ws.Range("H80").Value = 61996
$ws.Range("J80").Value = 61996
$ws.Range("L80").Value = 61996
$ws.Range("N80").Value = -63992

# Row 83: All's Fair in Highborn Assassination (L) | Titanium Hoplon
$ws.Range("H83").Value = 61996
$ws.Range("J83").Value = 61996
$ws.Range("L83").Value = 185988
$ws.Range("N83").Value = -195972

# Row 88: The Mast Chance | Adamantite Rivets
$ws.Range("H88").Value = 52452320
$ws.Range("I88").Value = 85717120
$ws.Range("J88").Value = 10115306
$ws.Range("K88").Value = 85717120
$ws.Range("L88").Value = 10115306
$ws.Range("M88").Value = -85716714
$ws.Range("N88").Value = -10116118

# Row 91: The Rose and the Riveter (L) | Adamantite Rivets
$ws.Range("H91").Value = 52452320
$ws.Range("I91").Value = 85717120
$ws.Range("J91").Value = 10115306
$ws.Range("K91").Value = 85717120
$ws.Range("L91").Value = 10115306
$ws.Range("M91").Value = -85715716
$ws.Range("N91").Value = -10118114

# Row 109: A Head of Demand | Deepgold Helm of Fending
$ws.Range("H109").Value = 37369
$ws.Range("J109").Value = 37369
$ws.Range("L109").Value = 37369
$ws.Range("N109").Value = -40143

# Row 117: Signed, Shield, Delivered | Titanbronze Tower Shield
$ws.Range("H117").Value = 46998
$ws.Range("J117").Value = 46998
$ws.Range("L117").Value = 46998
$ws.Range("N117").Value = -56176

# Row 118: A Budding Business | Titanbronze Headband of Scouting
$ws.Range("H118").Value = 49998
$ws.Range("J118").Value = 49998
$ws.Range("L118").Value = 49998
$ws.Range("N118").Value = -53312

# Row 120: One Foot Forward | Dwarven Mythril Shoes of Maiming
$ws.Range("H120").Value = 46104
$ws.Range("J120").Value = 46104
$ws.Range("L120").Value = 46104
$ws.Range("N120").Value = -55780

# Row 138: Don't Ask about the Rivets | Titanium Gold Helm of Casting
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

# ---------------- BSM ----------------
$ws = $wb.Worksheets.Item("BSM")

# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 2600
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 1800
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -4046

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 2600
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 1800
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -20232

# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 2146
$ws.Range("I105").Value = 2062.3684
$ws.Range("J105").Value = 2215.087
$ws.Range("K105").Value = 2062.3684
$ws.Range("L105").Value = 2215.087
$ws.Range("M105").Value = -315.3683999999998
$ws.Range("N105").Value = -5709.087

# Row 117: Idol Hands | Titanbronze Chakrams
$ws.Range("H117").Value = 49911.332
$ws.Range("J117").Value = 49911.332
$ws.Range("L117").Value = 49911.332
$ws.Range("N117").Value = -59089.332

# Row 119: Bae Blade | Dwarven Mythril Uchigatana
$ws.Range("H119").Value = 46711
$ws.Range("J119").Value = 46711
$ws.Range("L119").Value = 46711
$ws.Range("N119").Value = -56387

# Row 120: Under the Fool Moon | Dwarven Mythril Pistol
$ws.Range("H120").Value = 45761
$ws.Range("J120").Value = 45761
$ws.Range("L120").Value = 45761
$ws.Range("N120").Value = -55437

# Row 130: Annals of the Empire I | Chondrite Magitek Axe
$ws.Range("H130").Value = 41336
$ws.Range("J130").Value = 41336
$ws.Range("L130").Value = 41336
$ws.Range("N130").Value = -51376

# ---------------- CRP ----------------
$ws = $wb.Worksheets.Item("CRP")

# Row 105: Zelkova, My Love | Zelkova Lumber
$ws.Range("H105").Value = 2145.1365
$ws.Range("I105").Value = 2059.65
$ws.Range("K105").Value = 2059.65
$ws.Range("M105").Value = -312.6500000000001

# Row 116: The Right Tool for the Job | Sandteak Rod
$ws.Range("H116").Value = 59984.75
$ws.Range("J116").Value = 59984.75
$ws.Range("L116").Value = 59984.75
$ws.Range("N116").Value = -69162.75

# Row 118: A Miss and a Hit | Sandteak Longbow
$ws.Range("H118").Value = 44716
$ws.Range("J118").Value = 44716
$ws.Range("L118").Value = 44716
$ws.Range("N118").Value = -48030

# Row 139: Weaving a Path | Acacia Spinning Wheel
$ws.Range("H139").Value = 64052.715
$ws.Range("I139").Value = 42000
$ws.Range("J139").Value = 67728.164
$ws.Range("K139").Value = 42000
$ws.Range("L139").Value = 67728.164
$ws.Range("M139").Value = -36860
$ws.Range("N139").Value = -78008.164

# ---------------- CUL ----------------
$ws = $wb.Worksheets.Item("CUL")

# Row 22: A Total Nut Job | Walnut Bread
$ws.Range("H22").Value = 26766.5
$ws.Range("I22").Value = 18772
$ws.Range("J22").Value = 50750
$ws.Range("K22").Value = 56316
$ws.Range("L22").Value = 152250
$ws.Range("M22").Value = -56147
$ws.Range("N22").Value = -152588

# Row 27: Brain Food | Walnut Bread
$ws.Range("H27").Value = 26766.5
$ws.Range("I27").Value = 18772
$ws.Range("J27").Value = 50750
$ws.Range("K27").Value = 56316
$ws.Range("L27").Value = 152250
$ws.Range("M27").Value = -56214
$ws.Range("N27").Value = -152454

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 26406.518
$ws.Range("I131").Value = 9499.362999999999
$ws.Range("J131").Value = 28989.555
$ws.Range("K131").Value = 28498.089
$ws.Range("L131").Value = 86968.66500000001
$ws.Range("M131").Value = -23458.089
$ws.Range("N131").Value = -97048.66500000001

# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 1558.8438
$ws.Range("I132").Value = 957.1
$ws.Range("J132").Value = 1832.3636
$ws.Range("K132").Value = 8613.9
$ws.Range("L132").Value = 16491.2724
$ws.Range("M132").Value = -6083.9
$ws.Range("N132").Value = -21551.2724

# Row 133: Friends Are Food | Boiled Alpaca Steak
$ws.Range("H133").Value = 4291.25
$ws.Range("I133").Value = 4291.25
$ws.Range("K133").Value = 12873.75
$ws.Range("M133").Value = -7813.75

# Row 134: Don't Knock It Till You've Tried It | Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 62752356
$ws.Range("I134").Value = 66935444
$ws.Range("J134").Value = 6033
$ws.Range("K134").Value = 200806332
$ws.Range("L134").Value = 18099
$ws.Range("M134").Value = -200801262
$ws.Range("N134").Value = -28239

# Row 137: Creative Chocolate | Gateau au Chocolat
$ws.Range("H137").Value = 47627416
$ws.Range("I137").Value = 3328.1667
$ws.Range("J137").Value = 111126200
$ws.Range("K137").Value = 9984.500100000001
$ws.Range("L137").Value = 333378600
$ws.Range("M137").Value = -4884.500100000001
$ws.Range("N137").Value = -333388800

# Row 139: Najoothie | Wild Banana Blend
$ws.Range("H139").Value = 9458.223
$ws.Range("I139").Value = 14624
$ws.Range("J139").Value = 3001
$ws.Range("K139").Value = 43872
$ws.Range("L139").Value = 9003
$ws.Range("M139").Value = -38732
$ws.Range("N139").Value = -19283

# ---------------- GSM ----------------
$ws = $wb.Worksheets.Item("GSM")

# Row 110: Slimming Down | Stonegold Rapier
$ws.Range("H110").Value = 48702
$ws.Range("J110").Value = 48702
$ws.Range("L110").Value = 48702
$ws.Range("N110").Value = -56882

# ---------------- LTW ----------------
$ws = $wb.Worksheets.Item("LTW")

# Row 82: Trainin' the Neck | Dragon Leather
$ws.Range("H82").Value = 27780078
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

# Row 85: Training Is Only Skintight (L) | Dragon Leather
$ws.Range("H85").Value = 27780078
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# Row 111: Glove Me Tender | Gliderskin Gloves of Striking
$ws.Range("H111").Value = 44383
$ws.Range("J111").Value = 44383
$ws.Range("L111").Value = 44383
$ws.Range("N111").Value = -52563

# Row 121: A Shoe In | Swallowskin Shoes of Healing
$ws.Range("H121").Value = 44412
$ws.Range("J121").Value = 44412
$ws.Range("L121").Value = 44412
$ws.Range("N121").Value = -47906

# ---------------- WVR ----------------
$ws = $wb.Worksheets.Item("WVR")

# Row 16: Keep It under Wraps | Cotton Turban
$ws.Range("H16").Value = 44210
$ws.Range("J16").Value = 44210
$ws.Range("L16").Value = 44210
$ws.Range("N16").Value = -44794

# Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws.Range("H81").Value = 1855.579
$ws.Range("I81").Value = 1505.2941
$ws.Range("J81").Value = 4833
$ws.Range("K81").Value = 3010.5882
$ws.Range("L81").Value = 9666
$ws.Range("M81").Value = -1949.5882
$ws.Range("N81").Value = -11788

# Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws.Range("H84").Value = 1855.579
$ws.Range("I84").Value = 1505.2941
$ws.Range("J84").Value = 4833
$ws.Range("K84").Value = 15052.941
$ws.Range("L84").Value = 48330
$ws.Range("M84").Value = -9748.941000000001
$ws.Range("N84").Value = -58938

# Row 119: A Job Well Done | Dwarven Cotton Gaskins of Fending
$ws.Range("H119").Value = 43661.332
$ws.Range("J119").Value = 43661.332
$ws.Range("L119").Value = 43661.332
$ws.Range("N119").Value = -53337.332
